$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. Update the 10" pizza item names on Sheet1 (A5:A8) to say "... 10' Inch"
#    (order matches how the shared-string table ends up laid out on save)
$ws1.Range("A7").Value = "Pizza Onion and Capsicum 10' Inch"
$ws1.Range("A8").Value = "Pizza Paneer, Veggie ( Onion and Capsicum and corn) 10' Inch"
$ws1.Range("A6").Value = "Pizza Corn 10' Inch"
$ws1.Range("A5").Value = "Pizza Margarita 10' Inch"

# 2. Capture the four rows (Veg Thail, Veg Special Thail, Chole Bhature, Chole chawal)
#    that are being relocated from Sheet1 to Sheet2 before we delete them.
$moveRows = @(19, 20, 21, 22)
$captured = @()
foreach ($r in $moveRows) {
    $rowData = @{
        A = $ws1.Range("A$r").Value2
        C = $ws1.Range("C$r").Value2
        D = $ws1.Range("D$r").Value2
    }
    $captured += $rowData
}

# Grab the currency number format used for prices so the new rows match.
$priceFormat = $ws1.Range("C2").NumberFormat

# 3. Remove those four rows from Sheet1 (rows below shift up automatically).
$ws1.Rows("19:22").Delete() | Out-Null

# 4. Append the captured rows onto the end of Sheet2 (new rows 8-11).
$destRow = $ws2.UsedRange.Rows.Count + 1
foreach ($rowData in $captured) {
    $ws2.Range("A$destRow").Value = $rowData.A
    $ws2.Range("C$destRow").Value = $rowData.C
    $ws2.Range("C$destRow").NumberFormat = $priceFormat
    $ws2.Range("D$destRow").Value = $rowData.D
    $destRow = $destRow + 1
}

# 5. Restore the on-screen selections recorded in the workbook (Sheet2 first,
#    then Sheet1 last so Sheet1 stays the active/tabSelected sheet).
$ws2.Range("A15").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A6").Select() | Out-Null
